# The workbook tracks weekly price observations for "Ciboulette" at the
# "Vega Modelo de Temuco" market. This edit adds one more weekly
# observation: a new row is inserted right after the existing row 74
# (pushing all subsequent data rows down by one), and the new row is
# populated with the same Volumen/Precio/Unidad/Origen values as the
# (now shifted) row that used to be row 74, but dated one week later
# (serial 44533 = 2021-12-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 74; rows 74-184 shift down to 75-185,
# and the new row 74 inherits formatting from the row above it (so D74
# keeps the date number format used throughout column D).
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(74, 1).Value = 10
$ws.Cells.Item(74, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(74, 3).Value = "La Araucanía"
$ws.Cells.Item(74, 4).Value = 44533
$ws.Cells.Item(74, 5).Value = 9
$ws.Cells.Item(74, 6).Value = 100112039
$ws.Cells.Item(74, 7).Value = "Ciboulette"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 65
$ws.Cells.Item(74, 11).Value = 5000
$ws.Cells.Item(74, 12).Value = 5000
$ws.Cells.Item(74, 13).Value = 5000
$ws.Cells.Item(74, 14).Value = "`$/docena de atados"
$ws.Cells.Item(74, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(74, 16).Value = 1667
$ws.Cells.Item(74, 17).Value = 3
$ws.Cells.Item(74, 18).Value = "Hortaliza"
